# This script applies a large set of individual cell-value corrections to
# Sheet1 of the "seminare_bez_seminaricich" workbook. The underlying change
# is a row-level re-shuffle of previously-swapped datum/hodinaSkutOd/
# hodinaSkutDo/idno/jmena values back to their correct rows (commit message:
# "Ignore jsem ve spatny vetvi" -- i.e. reverting edits made on the wrong
# branch). No rows/columns are added or removed and no formatting changes;
# only specific cell values change, so we set each one explicitly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '17.4.2024'
$ws.Range("E3").Value = '24.4.2024'
$ws.Range("E4").Value = '15.5.2024'
$ws.Range("E5").Value = '14.5.2024'
$ws.Range("G5").Value = '15:50'
$ws.Range("E6").Value = '7.5.2024'
$ws.Range("G6").Value = '16:50'
$ws.Range("E7").Value = '30.4.2024'
$ws.Range("G7").Value = '15:50'
$ws.Range("E9").Value = '23.4.2024'
$ws.Range("G9").Value = '16:50'
$ws.Range("B21").Value = 251
$ws.Range("L21").Value = 'Fišer Jiří, Mgr. Ph.D.'
$ws.Range("B22").Value = 2220
$ws.Range("L22").Value = 'Škvor Jiří, RNDr. Ph.D.'
$ws.Range("F48").Value = '17:00'
$ws.Range("G48").Value = '18:50'
$ws.Range("E51").Value = '15.5.2024'
$ws.Range("E52").Value = '24.4.2024'
$ws.Range("F53").Value = '13:00'
$ws.Range("G53").Value = '14:50'
$ws.Range("E62").Value = '30.4.2024'
$ws.Range("F62").Value = '08:00'
$ws.Range("G62").Value = '09:50'
$ws.Range("E63").Value = '10.5.2024'
$ws.Range("E64").Value = '7.5.2024'
$ws.Range("F64").Value = '08:00'
$ws.Range("G64").Value = '09:50'
$ws.Range("B65").Value = 3521
$ws.Range("E65").Value = '13.5.2024'
$ws.Range("F65").Value = '09:00'
$ws.Range("G65").Value = '10:50'
$ws.Range("L65").Value = 'Posel Zbyšek, doc. RNDr. Ph.D.'
$ws.Range("B66").Value = 3521
$ws.Range("E66").Value = '29.4.2024'
$ws.Range("F66").Value = '09:00'
$ws.Range("G66").Value = '10:50'
$ws.Range("L66").Value = 'Posel Zbyšek, doc. RNDr. Ph.D.'
$ws.Range("E68").Value = '6.5.2024'
$ws.Range("E69").Value = '17.5.2024'
$ws.Range("F69").Value = '10:00'
$ws.Range("G69").Value = '11:50'
$ws.Range("B70").Value = 2317
$ws.Range("E70").Value = '12.4.2024'
$ws.Range("F70").Value = '10:00'
$ws.Range("G70").Value = '11:50'
$ws.Range("L70").Value = 'Kuba Pavel, Ing. Ph.D.'
$ws.Range("B71").Value = 3521
$ws.Range("E71").Value = '22.4.2024'
$ws.Range("F71").Value = '09:00'
$ws.Range("G71").Value = '10:50'
$ws.Range("L71").Value = 'Posel Zbyšek, doc. RNDr. Ph.D.'
$ws.Range("B72").Value = 2317
$ws.Range("E72").Value = '14.5.2024'
$ws.Range("F72").Value = '08:00'
$ws.Range("G72").Value = '09:50'
$ws.Range("L72").Value = 'Kuba Pavel, Ing. Ph.D.'
$ws.Range("E73").Value = '19.4.2024'
$ws.Range("B74").Value = 3521
$ws.Range("E74").Value = '15.4.2024'
$ws.Range("F74").Value = '09:00'
$ws.Range("G74").Value = '10:50'
$ws.Range("L74").Value = 'Posel Zbyšek, doc. RNDr. Ph.D.'
$ws.Range("E75").Value = '26.4.2024'
$ws.Range("B76").Value = 2317
$ws.Range("E76").Value = '3.5.2024'
$ws.Range("F76").Value = '10:00'
$ws.Range("G76").Value = '11:50'
$ws.Range("L76").Value = 'Kuba Pavel, Ing. Ph.D.'
$ws.Range("B77").Value = 2317
$ws.Range("E77").Value = '16.4.2024'
$ws.Range("F77").Value = '08:00'
$ws.Range("G77").Value = '09:50'
$ws.Range("L77").Value = 'Kuba Pavel, Ing. Ph.D.'
$ws.Range("E81").Value = '10.5.2024'
$ws.Range("E83").Value = '26.4.2024'
$ws.Range("E84").Value = '23.4.2024'
$ws.Range("E85").Value = '30.4.2024'
$ws.Range("E86").Value = '14.5.2024'
$ws.Range("E88").Value = '7.5.2024'
$ws.Range("E92").Value = '12.4.2024'
$ws.Range("G92").Value = '12:50'
$ws.Range("B94").Value = 6259
$ws.Range("E94").Value = '3.5.2024'
$ws.Range("L94").Value = 'Škvára Jiří, RNDr. Ph.D.'
$ws.Range("E95").Value = ""
$ws.Range("F95").Value = '00:00'
$ws.Range("G95").Value = '00:00'
$ws.Range("B96").Value = 251
$ws.Range("L96").Value = 'Fišer Jiří, Mgr. Ph.D.'
$ws.Range("E97").Value = '19.4.2024'
$ws.Range("G97").Value = '13:50'
$ws.Range("E98").Value = '17.5.2024'
$ws.Range("G98").Value = '13:50'
$ws.Range("B99").Value = 251
$ws.Range("E99").Value = '12.4.2024'
$ws.Range("F99").Value = '12:00'
$ws.Range("G99").Value = '12:50'
$ws.Range("L99").Value = 'Fišer Jiří, Mgr. Ph.D.'
$ws.Range("B100").Value = 6259
$ws.Range("E100").Value = '26.4.2024'
$ws.Range("L100").Value = 'Škvára Jiří, RNDr. Ph.D.'
$ws.Range("B101").Value = 6973
$ws.Range("E101").Value = ""
$ws.Range("F101").Value = '00:00'
$ws.Range("G101").Value = '00:00'
$ws.Range("L101").Value = 'Beránek Pavel, Ing. Mgr.'
$ws.Range("E102").Value = '19.4.2024'
$ws.Range("F102").Value = '12:00'
$ws.Range("G102").Value = '13:50'
$ws.Range("B104").Value = 2776
$ws.Range("E104").Value = ""
$ws.Range("F104").Value = '00:00'
$ws.Range("G104").Value = '00:00'
$ws.Range("L104").Value = 'Krejčí Jan, RNDr. Ph.D.'
$ws.Range("B105").Value = 251
$ws.Range("E105").Value = '3.5.2024'
$ws.Range("F105").Value = '12:00'
$ws.Range("G105").Value = '13:50'
$ws.Range("L105").Value = 'Fišer Jiří, Mgr. Ph.D.'
$ws.Range("E106").Value = '24.4.2024'
$ws.Range("E107").Value = '15.5.2024'
$ws.Range("E109").Value = '18.4.2024'
$ws.Range("F109").Value = '09:00'
$ws.Range("G109").Value = '11:50'
$ws.Range("E110").Value = '2.5.2024'
$ws.Range("E111").Value = '3.5.2024'
$ws.Range("F111").Value = '16:00'
$ws.Range("G111").Value = '18:50'
$ws.Range("E112").Value = '19.4.2024'
$ws.Range("F112").Value = '16:00'
$ws.Range("G112").Value = '18:50'
$ws.Range("E113").Value = '9.5.2024'
$ws.Range("F113").Value = '09:00'
$ws.Range("G113").Value = '11:50'
$ws.Range("E114").Value = '16.5.2024'
$ws.Range("E115").Value = '11.4.2024'
$ws.Range("E116").Value = '25.4.2024'
$ws.Range("F116").Value = '09:00'
$ws.Range("G116").Value = '11:50'
$ws.Range("E117").Value = '17.5.2024'
$ws.Range("F117").Value = '16:00'
$ws.Range("G117").Value = '18:50'
